# Update the "想去人数" (people interested) counts that changed between
# the previous and newly generated gh-pages data snapshot.
# Both the "展览" sheet and the "全部类型" sheet carry the same rows,
# so the same updates must be applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 25
    "F5"  = 63
    "F6"  = 2330
    "F8"  = 1429
    "F12" = 380
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
